# Updated cryptos list with latest scrape data (price + 1h volume change).
# Numeric-looking price strings are written with a leading apostrophe so
# Excel keeps them as literal text (matching the original inline-string cells)
# instead of auto-converting them into numeric cell values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.700.91"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.475.84"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'320.22"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").Value = "'92.41"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'0.0863"
$ws.Range("E10").Value = "  +6.67%  "
$ws.Range("D11").Value = "'33.09"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "2.857.19"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "'6.91"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").Value = "2.464.95"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "'0.795"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "41.621.14"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "0.0₃0943"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'11.27"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").Value = "'239.72"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").Value = "'2.76"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'25.07"
$ws.Range("E27").Value = "  +3.25%  "
$ws.Range("D28").Value = "'2.24"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").Value = "'36.70"
$ws.Range("E30").Value = "  +4.43%  "
$ws.Range("D31").Value = "'157.53"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'0.0765"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "'17.16"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.116"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").Value = "'4.03"
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("D42").Value = "'2.44"
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("D43").Value = "1.997.48"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").Value = "'18.72"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("E47").Value = "  +5.59%  "
$ws.Range("D48").Value = "2.752.65"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").Value = "'97.63"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("D50").Value = "'76.08"
$ws.Range("E50").Value = "  +5.73%  "
$ws.Range("D51").Value = "'67.47"
$ws.Range("E51").Value = "  +0.95%  "
